$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing X5/Y5 values on the existing last row
$ws.Range("X5").Value = 0.11999500000000296
$ws.Range("Y5").Value = "Up"

# Append a new row (row 6) of scan results.
# Seed formatted cells (date + percentages) by copying the format from the
# row above so we reuse the existing style indexes instead of minting new ones.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("S5").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("T5").Copy()
$ws.Range("T6").PasteSpecial(-4122)

$ws.Range("A6").Value = 42647.887002314812
$ws.Range("B6").Value = -3
$ws.Range("C6").Value = "Neutral"
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 16039
$ws.Range("F6").Value = 956
$ws.Range("G6").Value = 63
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 74
$ws.Range("J6").Value = 24
$ws.Range("K6").Value = 32830
$ws.Range("L6").Value = 167
$ws.Range("M6").Value = 94
$ws.Range("N6").Value = 39
$ws.Range("O6").Value = 13
$ws.Range("P6").Value = "Named"
$ws.Range("Q6").Value = 42.459412013272512
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = -0.0112
$ws.Range("T6").Value = -0.036700000000000003
$ws.Range("U6").Value = 14.56
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = -2
